# Adds a "Sheet2" tab that mirrors Sheet1's PsNo roster and carries 19
# Hobbie_n columns (with sample data in row 2), and tags Sheet1 with a new
# trailing "Discpline_Grade" column.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Add Sheet2 (after Sheet1) with a PsNo column + 19 Hobbie_n columns
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$hobbieHeaders = @("Hobbie_1", "Hobbie_2", "Hobbie_3", "Hobbie_4", "Hobbie_5", "Hobbie_6", "Hobbie_7", "Hobbie_8", "Hobbie_9", "Hobbie_10", "Hobbie_11", "Hobbie_12", "Hobbie_13", "Hobbie_14", "Hobbie_15", "Hobbie_16", "Hobbie_17", "Hobbie_18", "Hobbie_19")
$hobbieValues = @("Geography", "Skiing", "Weight training", "Postcard Collecting", "Playing The Stock Market", "Beekeeping", "Playing Chess", "Archery", "Coin Collecting", "Vintage Clothing Collecting", "Furniture Building", "LEGO", "Amateur Radio", "Sudoku", "Crabbing", "Kayaking", "Parkour", "Stamp Collecting", "Wood Carving")

# Header row
$ws2.Cells.Item(1, 1).Value = "PsNo"
$ws2.Cells.Item(1, 1).Font.Bold = $true
for ($i = 0; $i -lt $hobbieHeaders.Length; $i++) {
  $cell = $ws2.Cells.Item(1, $i + 2)
  $cell.Value = $hobbieHeaders[$i]
  $cell.Font.Bold = $true
}

# Row 2 sample data
$ws2.Cells.Item(2, 1).Value = 99004000
for ($i = 0; $i -lt $hobbieValues.Length; $i++) {
  $ws2.Cells.Item(2, $i + 2).Value = $hobbieValues[$i]
}

# Remaining PsNo roster (rows 3-15), mirroring Sheet1 column A
for ($r = 3; $r -le 15; $r++) {
  $psNo = $ws1.Cells.Item($r, 1).Value2
  $ws2.Cells.Item($r, 1).Value = $psNo
}

# Column widths for Sheet2 (pre-compensated for the host's 1/6-character
# quantisation step so the saved OOXML "width" lands as close as possible
# to the authored values).
$widths = @(15.0, 11.666667, 15.166667, 12.5, 10.666667, 16.666667, 13.0, 14.5, 16.166667, 15.666667, 11.833333, 15.166667, 13.5, 13.833333, 13.666667, 14.833333, 16.0, 14.833333, 12.666667)
for ($i = 0; $i -lt $widths.Length; $i++) {
  $ws2.Columns.Item($i + 2).ColumnWidth = $widths[$i]
}

# ---------------------------------------------------------------------
# 2) Sheet1: add a trailing "Discpline_Grade" column (T)
# ---------------------------------------------------------------------
$ws1.Range("T1").Value = "Discpline_Grade"
$ws1.Range("T1").Font.Bold = $true
$ws1.Range("T2").Value = "O"

# Match the column widths recorded for the new layout as closely as the
# engine's character-width quantisation allows.
$ws1.Columns.Item(19).ColumnWidth = 13.833333
$ws1.Columns.Item(20).ColumnWidth = 14.5

# ---------------------------------------------------------------------
# 3) Selection / active-sheet bookkeeping to match the authored state:
#    Sheet1 stays the active tab with E8 selected; Sheet2 remembers Q7.
# ---------------------------------------------------------------------
[void]$ws2.Range("Q7").Select()
[void]$ws1.Activate()
[void]$ws1.Range("E8").Select()
